$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, shifting existing rows 34..149 down to 35..150
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new data record
$ws.Cells.Item(34, 1).Value = 3
$ws.Cells.Item(34, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(34, 3).Value = "Coquimbo"
$ws.Cells.Item(34, 4).Value = 44620
$ws.Cells.Item(34, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(34, 5).Value = 5
$ws.Cells.Item(34, 6).Value = 100112052
$ws.Cells.Item(34, 7).Value = "Albahaca"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 110
$ws.Cells.Item(34, 11).Value = 5500
$ws.Cells.Item(34, 12).Value = 6000
$ws.Cells.Item(34, 13).Value = 5727
$ws.Cells.Item(34, 14).Value = "`$/docena de matas"
$ws.Cells.Item(34, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(34, 16).Value = 954
$ws.Cells.Item(34, 17).Value = 6
$ws.Cells.Item(34, 18).Value = "Hortaliza"
